$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("H3").Value = 3.3
$ws.Range("L3").Value = 1.5
$ws.Range("M3").Value = 2.5
$ws.Range("N3").Value = 2.5
$ws.Range("O3").Value = 1.5
$ws.Range("R3").Value = 2.25
$ws.Range("S3").Value = 1.57
$ws.Range("V3").Value = 9.5
$ws.Range("X3").Value = 19
# Row 5
$ws.Range("G5").Value = 1.75
$ws.Range("I5").Value = 4.75
$ws.Range("J5").Value = 1.08
$ws.Range("K5").Value = 8
$ws.Range("Z5").Value = 8
$ws.Range("AA5").Value = 6.5
$ws.Range("AE5").Value = 11
# Row 6
$ws.Range("G6").Value = 1.7
$ws.Range("H6").Value = 3.7
$ws.Range("J6").Value = 1.05
$ws.Range("K6").Value = 11
$ws.Range("L6").Value = 1.29
$ws.Range("M6").Value = 3.5
$ws.Range("N6").Value = 1.93
$ws.Range("O6").Value = 1.93
$ws.Range("P6").Value = 1.4
$ws.Range("Q6").Value = 2.75
$ws.Range("R6").Value = 1.83
$ws.Range("S6").Value = 1.83
$ws.Range("T6").Value = 7
$ws.Range("U6").Value = 8
$ws.Range("X6").Value = 13
$ws.Range("Y6").Value = 26
$ws.Range("Z6").Value = 10
$ws.Range("AA6").Value = 7
$ws.Range("AB6").Value = 15
$ws.Range("AD6").Value = 251
$ws.Range("AE6").Value = 13
$ws.Range("AF6").Value = 26
$ws.Range("AG6").Value = 15
# Row 9
$ws.Range("I9").Value = 9.5
$ws.Range("K9").Value = 13
$ws.Range("AB9").Value = 21
$ws.Range("AD9").Value = 1250
$ws.Range("AG9").Value = 29
$ws.Range("AH9").Value = 126
# Row 10
$ws.Range("G10").Value = 1.44
$ws.Range("H10").Value = 3.9
$ws.Range("I10").Value = 9
$ws.Range("J10").Value = 1.1
$ws.Range("K10").Value = 7
$ws.Range("L10").Value = 1.4
$ws.Range("M10").Value = 2.75
$ws.Range("N10").Value = 2.3
$ws.Range("O10").Value = 1.6
$ws.Range("T10").Value = 5
$ws.Range("AI10").Value = 81
# Row 11
$ws.Range("J11").Value = 1.08
$ws.Range("K11").Value = 8
# Row 14
$ws.Range("G14").Value = 2.22
$ws.Range("H14").Value = 3.3
$ws.Range("I14").Value = 2.95
$ws.Range("R14").Value = 1.7
$ws.Range("S14").Value = 1.91
$ws.Range("T14").Value = 7.7
$ws.Range("U14").Value = 10.75
$ws.Range("V14").Value = 9
$ws.Range("W14").Value = 22
$ws.Range("X14").Value = 18.5
$ws.Range("Y14").Value = 29
$ws.Range("AA14").Value = 6.4
$ws.Range("AB14").Value = 14
$ws.Range("AE14").Value = 9.25
$ws.Range("AF14").Value = 15.5
$ws.Range("AG14").Value = 10.5
$ws.Range("AH14").Value = 37
$ws.Range("AI14").Value = 25
$ws.Range("AJ14").Value = 32
# Row 15
$ws.Range("G15").Value = 4.05
$ws.Range("H15").Value = 3.75
$ws.Range("I15").Value = 1.72
$ws.Range("N15").Value = 1.5
$ws.Range("O15").Value = 2.25
$ws.Range("R15").Value = 1.5
$ws.Range("S15").Value = 2.25
$ws.Range("T15").Value = 17
$ws.Range("U15").Value = 28
$ws.Range("V15").Value = 13.5
$ws.Range("W15").Value = 65
$ws.Range("X15").Value = 32
$ws.Range("Y15").Value = 30
$ws.Range("Z15").Value = 15.5
$ws.Range("AA15").Value = 7.8
$ws.Range("AB15").Value = 12
$ws.Range("AC15").Value = 40
$ws.Range("AD15").Value = 250
$ws.Range("AE15").Value = 9.75
$ws.Range("AJ15").Value = 18.5
# Row 19
$ws.Range("G19").Value = 2.05
$ws.Range("I19").Value = 3.25
$ws.Range("J19").Value = 1.04
$ws.Range("K19").Value = 9
$ws.Range("L19").Value = 1.22
$ws.Range("M19").Value = 4
$ws.Range("N19").Value = 1.75
$ws.Range("O19").Value = 2.05
$ws.Range("T19").Value = 9
$ws.Range("Z19").Value = 12
$ws.Range("AA19").Value = 7
$ws.Range("AG19").Value = 12
# Row 21
$ws.Range("G21").Value = 1.57
$ws.Range("H21").Value = 4
$ws.Range("I21").Value = 5
$ws.Range("N21").Value = 1.57
$ws.Range("O21").Value = 2.35
$ws.Range("P21").Value = 1.29
$ws.Range("Q21").Value = 3.5
$ws.Range("R21").Value = 1.67
$ws.Range("S21").Value = 2.1
$ws.Range("T21").Value = 9
$ws.Range("U21").Value = 9
$ws.Range("AE21").Value = 17
$ws.Range("AI21").Value = 34
# Row 22
$ws.Range("G22").Value = 3.6
$ws.Range("H22").Value = 4
$ws.Range("I22").Value = 1.83
$ws.Range("J22").Value = 1.02
$ws.Range("K22").Value = 12
$ws.Range("L22").Value = 1.11
$ws.Range("M22").Value = 6
$ws.Range("N22").Value = 1.4
$ws.Range("O22").Value = 2.75
$ws.Range("T22").Value = 19
$ws.Range("U22").Value = 23
$ws.Range("V22").Value = 13
$ws.Range("X22").Value = 23
$ws.Range("Z22").Value = 21
$ws.Range("AA22").Value = 8.5
# Row 23
$ws.Range("G23").Value = 2.27
$ws.Range("H23").Value = 2.92
$ws.Range("I23").Value = 3.25
$ws.Range("L23").Value = 1.47
$ws.Range("M23").Value = 2.35
$ws.Range("N23").Value = 2.32
$ws.Range("O23").Value = 1.47
$ws.Range("P23").Value = 1.5
$ws.Range("Q23").Value = 2.25
$ws.Range("T23").Value = 5.9
$ws.Range("V23").Value = 9.75
$ws.Range("W23").Value = 23
$ws.Range("X23").Value = 23
$ws.Range("Z23").Value = 6.6
$ws.Range("AA23").Value = 5.8
$ws.Range("AB23").Value = 17.5
$ws.Range("AC23").Value = 110
$ws.Range("AE23").Value = 7.7
$ws.Range("AF23").Value = 15.5
$ws.Range("AH23").Value = 45
$ws.Range("AI23").Value = 35
